$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.290.30"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.914.53"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'488.60"
$ws.Range("E5").Value = "  +3.88%  "
$ws.Range("D6").Value = "'146.82"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "'43.06"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'10.83"
$ws.Range("E13").Value = "  +4.51%  "
$ws.Range("D14").Value = "4.539.00"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "3.905.18"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").Value = "'14.20"
$ws.Range("E16").Value = "  -6.04%  "
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "'19.88"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "68.372.48"
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("D21").Value = "'434.33"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "'14.92"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "'87.78"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'11.46"
$ws.Range("D26").Value = "'11.25"
$ws.Range("E26").Value = "  +10.72%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'38.19"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'5.73"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "'725.98"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("D31").Value = "'13.77"
$ws.Range("E31").Value = "  +2.66%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("E34").Value = "  +17.66%  "
$ws.Range("D35").Value = "'41.67"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "0.0₃0872"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("D37").Value = "'60.69"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("D38").Value = "'0.406"
$ws.Range("E38").Value = "  +21.04%  "
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +17.17%  "
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").Value = "'2.91"
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'3.33"
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "0.0₆0348"
$ws.Range("E50").Value = "  +34.67%  "
$ws.Range("D51").Value = "'144.90"
$ws.Range("E51").Value = "  -1.74%  "
